# Auto-generated script to apply scheduled-runner market data updates
$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsALC.Range("H96").Value = 761.6667
$wsALC.Range("I96").Value = 620.2222
$wsALC.Range("K96").Value = 1860.6666
$wsALC.Range("M96").Value = -487.6666
$wsALC.Range("H97").Value = 225
$wsALC.Range("J97").Value = 206.25
$wsALC.Range("L97").Value = 618.75
$wsALC.Range("N97").Value = -1610.75
$wsALC.Range("H137").Value = 38892.594
$wsALC.Range("I137").Value = 1570.6666
$wsALC.Range("J137").Value = 113536.445
$wsALC.Range("K137").Value = 4711.9998
$wsALC.Range("L137").Value = 340609.335
$wsALC.Range("M137").Value = -2161.9998
$wsALC.Range("N137").Value = -345709.335
$wsALC.Range("H138").Value = 3130.6445
$wsALC.Range("J138").Value = 4914.885
$wsALC.Range("L138").Value = 14744.655
$wsALC.Range("N138").Value = -25024.655

$wsARM = $wb.Worksheets.Item("ARM")
$wsARM.Range("H2").Value = 2105.5
$wsARM.Range("I2").Value = 2440.1
$wsARM.Range("J2").Value = 1547.8334
$wsARM.Range("K2").Value = 2440.1
$wsARM.Range("L2").Value = 1547.8334
$wsARM.Range("M2").Value = -2327.1
$wsARM.Range("N2").Value = -1773.8334
$wsARM.Range("H5").Value = 390
$wsARM.Range("I5").Value = 186.66667
$wsARM.Range("K5").Value = 186.66667
$wsARM.Range("M5").Value = -74.66667000000001
$wsARM.Range("H32").Value = 21280564
$wsARM.Range("I32").Value = 23258062
$wsARM.Range("J32").Value = 22456.75
$wsARM.Range("K32").Value = 23258062
$wsARM.Range("L32").Value = 22456.75
$wsARM.Range("M32").Value = -23257775
$wsARM.Range("N32").Value = -23030.75
$wsARM.Range("H61").Value = 3181.3635
$wsARM.Range("I61").Value = 2746.923
$wsARM.Range("J61").Value = 3808.889
$wsARM.Range("K61").Value = 2746.923
$wsARM.Range("L61").Value = 3808.889
$wsARM.Range("M61").Value = -2534.923
$wsARM.Range("N61").Value = -4232.889
$wsARM.Range("H110").Value = 3695.7368
$wsARM.Range("I110").Value = 3171.8572
$wsARM.Range("J110").Value = 5162.6
$wsARM.Range("K110").Value = 3171.8572
$wsARM.Range("L110").Value = 5162.6
$wsARM.Range("M110").Value = -1126.8572
$wsARM.Range("N110").Value = -9252.6
$wsARM.Range("H116").Value = 2105.5
$wsARM.Range("I116").Value = 2440.1
$wsARM.Range("J116").Value = 1547.8334
$wsARM.Range("K116").Value = 2440.1
$wsARM.Range("L116").Value = 1547.8334
$wsARM.Range("M116").Value = -146.0999999999999
$wsARM.Range("N116").Value = -6135.8334
$wsARM.Range("H122").Value = 1046271.44
$wsARM.Range("I122").Value = 1150675.2
$wsARM.Range("J122").Value = 2233
$wsARM.Range("K122").Value = 3452025.6
$wsARM.Range("L122").Value = 6699
$wsARM.Range("M122").Value = -3449575.6
$wsARM.Range("N122").Value = -11599
$wsARM.Range("H132").Value = 1756.5853
$wsARM.Range("I132").Value = 1448.862
$wsARM.Range("J132").Value = 2500.25
$wsARM.Range("K132").Value = 4346.586
$wsARM.Range("L132").Value = 7500.75
$wsARM.Range("M132").Value = -1816.586
$wsARM.Range("N132").Value = -12560.75
$wsARM.Range("H136").Value = 3181.3635
$wsARM.Range("I136").Value = 2746.923
$wsARM.Range("J136").Value = 3808.889
$wsARM.Range("K136").Value = 8240.769
$wsARM.Range("L136").Value = 11426.667
$wsARM.Range("M136").Value = -5690.769
$wsARM.Range("N136").Value = -16526.667

$wsBSM = $wb.Worksheets.Item("BSM")
$wsBSM.Range("H3").Value = 2105.5
$wsBSM.Range("I3").Value = 2440.1
$wsBSM.Range("J3").Value = 1547.8334
$wsBSM.Range("K3").Value = 2440.1
$wsBSM.Range("L3").Value = 1547.8334
$wsBSM.Range("M3").Value = -2326.1
$wsBSM.Range("N3").Value = -1775.8334
$wsBSM.Range("H4").Value = 390
$wsBSM.Range("I4").Value = 186.66667
$wsBSM.Range("K4").Value = 186.66667
$wsBSM.Range("M4").Value = -71.66667000000001
$wsBSM.Range("H134").Value = 1976.742
$wsBSM.Range("I134").Value = 1923.4138
$wsBSM.Range("K134").Value = 5770.2414
$wsBSM.Range("M134").Value = -3235.2414

$wsCRP = $wb.Worksheets.Item("CRP")
$wsCRP.Range("H16").Value = 2571
$wsCRP.Range("I16").Value = 1498.1666
$wsCRP.Range("J16").Value = 4716.6665
$wsCRP.Range("K16").Value = 1498.1666
$wsCRP.Range("L16").Value = 4716.6665
$wsCRP.Range("M16").Value = -1211.1666
$wsCRP.Range("N16").Value = -5290.6665
$wsCRP.Range("H31").Value = 6183.174
$wsCRP.Range("I31").Value = 6595.4287
$wsCRP.Range("J31").Value = 5541.8887
$wsCRP.Range("K31").Value = 6595.4287
$wsCRP.Range("L31").Value = 5541.8887
$wsCRP.Range("M31").Value = -6300.4287
$wsCRP.Range("N31").Value = -6131.8887
$wsCRP.Range("H34").Value = 6183.174
$wsCRP.Range("I34").Value = 6595.4287
$wsCRP.Range("J34").Value = 5541.8887
$wsCRP.Range("K34").Value = 6595.4287
$wsCRP.Range("L34").Value = 5541.8887
$wsCRP.Range("M34").Value = -6393.4287
$wsCRP.Range("N34").Value = -5945.8887
$wsCRP.Range("H99").Value = 2771
$wsCRP.Range("I99").Value = 2580.6667
$wsCRP.Range("J99").Value = 2999.4
$wsCRP.Range("K99").Value = 2580.6667
$wsCRP.Range("L99").Value = 2999.4
$wsCRP.Range("M99").Value = -1082.6667
$wsCRP.Range("N99").Value = -5995.4
$wsCRP.Range("H113").Value = 2571
$wsCRP.Range("I113").Value = 1498.1666
$wsCRP.Range("J113").Value = 4716.6665
$wsCRP.Range("K113").Value = 1498.1666
$wsCRP.Range("L113").Value = 4716.6665
$wsCRP.Range("M113").Value = 671.8334
$wsCRP.Range("N113").Value = -9056.666499999999
$wsCRP.Range("H126").Value = 2771
$wsCRP.Range("I126").Value = 2580.6667
$wsCRP.Range("J126").Value = 2999.4
$wsCRP.Range("K126").Value = 7742.000100000001
$wsCRP.Range("L126").Value = 8998.200000000001
$wsCRP.Range("M126").Value = -5272.000100000001
$wsCRP.Range("N126").Value = -13938.2
$wsCRP.Range("H132").Value = 1551.8889
$wsCRP.Range("I132").Value = 971.4828
$wsCRP.Range("J132").Value = 3956.4285
$wsCRP.Range("K132").Value = 2914.4484
$wsCRP.Range("L132").Value = 11869.2855
$wsCRP.Range("M132").Value = -384.4484000000002
$wsCRP.Range("N132").Value = -16929.2855

$wsCUL = $wb.Worksheets.Item("CUL")
$wsCUL.Range("H97").Value = 1012.38464
$wsCUL.Range("I97").Value = 597.75
$wsCUL.Range("K97").Value = 1793.25
$wsCUL.Range("M97").Value = -1297.25
$wsCUL.Range("H113").Value = 1815417.6
$wsCUL.Range("I113").Value = 6896972
$wsCUL.Range("J113").Value = 576.7143
$wsCUL.Range("K113").Value = 20690916
$wsCUL.Range("L113").Value = 1730.1429
$wsCUL.Range("M113").Value = -20688746
$wsCUL.Range("N113").Value = -6070.1429
$wsCUL.Range("H131").Value = 786.7931
$wsCUL.Range("I131").Value = 377
$wsCUL.Range("J131").Value = 1002.4737
$wsCUL.Range("K131").Value = 1131
$wsCUL.Range("L131").Value = 3007.4211
$wsCUL.Range("M131").Value = 3909
$wsCUL.Range("N131").Value = -13087.4211
$wsCUL.Range("H132").Value = 225361.05
$wsCUL.Range("I132").Value = 695.2121
$wsCUL.Range("J132").Value = 843192.0600000001
$wsCUL.Range("K132").Value = 6256.908899999999
$wsCUL.Range("L132").Value = 7588728.540000001
$wsCUL.Range("M132").Value = -3726.908899999999
$wsCUL.Range("N132").Value = -7593788.540000001
$wsCUL.Range("H137").Value = 12976.615
$wsCUL.Range("I137").Value = 4138.3335
$wsCUL.Range("J137").Value = 20552.285
$wsCUL.Range("K137").Value = 12415.0005
$wsCUL.Range("L137").Value = 61656.855
$wsCUL.Range("M137").Value = -7315.000499999998
$wsCUL.Range("N137").Value = -71856.855

$wsGSM = $wb.Worksheets.Item("GSM")
$wsGSM.Range("H102").Value = 1610.9062
$wsGSM.Range("I102").Value = 1680.625
$wsGSM.Range("J102").Value = 1401.75
$wsGSM.Range("K102").Value = 1680.625
$wsGSM.Range("L102").Value = 1401.75
$wsGSM.Range("M102").Value = -58.625
$wsGSM.Range("N102").Value = -4645.75
$wsGSM.Range("H132").Value = 3655.2942
$wsGSM.Range("I132").Value = 3100.125
$wsGSM.Range("J132").Value = 4148.778
$wsGSM.Range("K132").Value = 9300.375
$wsGSM.Range("L132").Value = 12446.334
$wsGSM.Range("M132").Value = -6770.375
$wsGSM.Range("N132").Value = -17506.334
